$d = $word.ActiveDocument

# The document has Track Changes on by default; temporarily disable it so
# the edits below land as plain content (as in the target) instead of
# tracked insertions/deletions. We restore the original setting at the end.
$origTrack = $d.TrackRevisions
$d.TrackRevisions = $false

# 1) Remove the old "_GoBack" bookmark from its current location (an empty
#    paragraph near the top of the "DATOS DEL RECURSO" block). It gets
#    re-created below, right before "El Memorial de Agravios" in the second
#    occurrence of the heading text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) First occurrence: " Refuerza tu aprendizaje: El Memorial de Agravios"
#    simply loses the "Refuerza tu aprendizaje: " prefix, becoming
#    " El Memorial de Agravios" (leading space kept, single run). Replace
#    only the first match in the document (the "Título del recurso" block).
$rng1 = $d.Content
[void]$rng1.Find.Execute("Refuerza tu aprendizaje: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1)

# 3) Second occurrence: locate the remaining "Refuerza tu aprendizaje: "
#    prefix (the one preceding "El Memorial de Agravios" in the "Título del
#    ejercicio" block).
$prefixRng = $d.Content
[void]$prefixRng.Find.Execute("Refuerza tu aprendizaje: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $prefixRng.Find.Found) {
    throw "Could not find the remaining 'Refuerza tu aprendizaje: ' prefix"
}
$prefixStart = $prefixRng.Start
$prefixEnd = $prefixRng.End

# Insert the _GoBack bookmark exactly where the prefix ends (i.e. right
# before "El Memorial de Agravios"); once the prefix text below is deleted,
# this sits precisely between the retained leading space and the title,
# splitting that text into two runs with the bookmark between them.
$bmRange = $d.Range($prefixEnd, $prefixEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Delete the "Refuerza tu aprendizaje: " prefix text itself, leaving
# " " + <bookmark> + "El Memorial de Agravios".
$delRange = $d.Range($prefixStart, $prefixEnd)
$delRange.Delete()

# Restore the document's original Track Changes setting.
$d.TrackRevisions = $origTrack
